# Update "想去人数" (attendee count) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 671
    3  = 510
    7  = 45
    8  = 2773
    9  = 4190
    10 = 104
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
